$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.374.59"
$ws.Range("E2").Value = "  -6.07%  "
$ws.Range("D3").Value = "3.187.41"
$ws.Range("E3").Value = "  -9.57%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.17"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -6.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.99"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -14.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "3.176.22"
$ws.Range("E8").Value = "  -9.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.538"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -11.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.168"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -14.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.30"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -11.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.491"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -16.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.09"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -18.34%  "
$ws.Range("E14").Value = "  -13.86%  "
$ws.Range("D15").Value = "3.681.53"
$ws.Range("E15").Value = "  -9.91%  "
$ws.Range("D16").Value = "66.304.32"
$ws.Range("E16").Value = "  -6.05%  "
$ws.Range("D17").Value = "3.180.86"
$ws.Range("E17").Value = "  -9.45%  "
$ws.Range("E18").Value = "  -6.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "529.20"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -14.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.06"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -16.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.92"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -16.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.749"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -15.48%  "
$ws.Range("E23").Value = "  -14.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.49"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -14.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.23"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -15.98%  "
$ws.Range("E27").Value = "  -18.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.15"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -17.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.96"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -13.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.78"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -15.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.54"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -15.79%  "
$ws.Range("E32").Value = "  -14.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.41"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -21.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "521.84"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -15.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.59"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -18.66%  "
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.93"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -7.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0848"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -15.79%  "
$ws.Range("E39").Value = "  -16.95%  "
$ws.Range("E40").Value = "  -17.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.124"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -14.81%  "
$ws.Range("D42").Value = "2.874.28"
$ws.Range("E42").Value = "  -14.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.60"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -26.07%  "
$ws.Range("E44").Value = "  -17.65%  "
$ws.Range("D45").Value = "0.0₃0577"
$ws.Range("E45").Value = "  -22.43%  "
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.77"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -20.31%  "
$ws.Range("E48").Value = "  -21.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.07"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -19.90%  "
$ws.Range("E50").Value = "  -14.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "122.03"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -8.74%  "
